$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 13 - "Ace KO e coli only"
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Ace KO e coli only"
$ws.Range("B13").Value = "Acetate KO Ecoli only 10-21-22"
$ws.Range("C13").Value = "https://anl.box.com/s/8lzgu1nf5xaduhldkfafwjpzcwtrhrxt"
$ws.Hyperlinks.Add($ws.Range("C13"), "https://anl.box.com/s/8lzgu1nf5xaduhldkfafwjpzcwtrhrxt") | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

$ws.Range("D13").Value = "10/21/2022"
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

$ws.Range("E13").Value = "EC Ace KO cells"
$ws.Range("H13").Value = "Acetate at various concentrations"
$ws.Range("I13").Value = "Maltose with acetate at various concentrations"

# ---------------------------------------------------------------------------
# Row 14 - "Metabolomics"
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Metabolomics"
$ws.Range("B14").Value = "8-19-22 Metabolomics TMS Panel EC_PF_acetate"
$ws.Range("C14").Value = "https://anl.box.com/s/90cvb2ddpxo7kv3ek91rx1vyceko8vb9"
$ws.Hyperlinks.Add($ws.Range("C14"), "https://anl.box.com/s/90cvb2ddpxo7kv3ek91rx1vyceko8vb9") | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null

$ws.Range("D14").Value = "8/19/2022"
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null

$ws.Range("E14").Value = "E. coli MG1655"
$ws.Range("F14").Value = "E. coli MG1655 with P. fluorescens"
$ws.Range("H14").Value = "Maltose alone"
$ws.Range("I14").Value = "Maltose and 4HB"
$ws.Range("J14").Value = "Maltose and acetate"

# ---------------------------------------------------------------------------
# Row 15 - "ODs for metabolomics"
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "ODs for metabolomics"
$ws.Range("B15").Value = "Appx. Cell ODs for Metabolomics 8-19-22 experiment"
$ws.Range("C15").Value = "https://anl.box.com/s/s1he202yyn8tye4i104e3nc0xaoxq5a3"
$ws.Hyperlinks.Add($ws.Range("C15"), "https://anl.box.com/s/s1he202yyn8tye4i104e3nc0xaoxq5a3") | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null

$ws.Range("D15").Value = "8/19/2022"
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null

$ws.Range("E15").Value = "E. coli MG1655"
$ws.Range("F15").Value = "E. coli MG1655 with P. fluorescens"
$ws.Range("H15").Value = "Maltose alone"
$ws.Range("I15").Value = "Maltose and 4HB"
$ws.Range("J15").Value = "Maltose and acetate"

# ---------------------------------------------------------------------------
# Column widths (approximate best-fit values, engine granularity limited)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.08
$ws.Columns.Item(2).ColumnWidth = 48.1
$ws.Columns.Item(9).ColumnWidth = 42.5

# ---------------------------------------------------------------------------
# Sheet view / selection
# ---------------------------------------------------------------------------
$ws.Range("C20").Select() | Out-Null
